# Add the new "edRVFL" model column (K) and update refreshed metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1 - copy style from existing header (e.g. J1) then set value.
$ws.Range("K1").Value = "edRVFL"
$ws.Range("K1").Style = $ws.Range("J1").Style

# New best_params cell K2 for edRVFL.
$ws.Range("K2").Value = "{'activation': 'relu', 'b_random_vec_range': [0, 10], 'lam': 1, 'n_layer': 16, 'n_nodes': 256, 'random_seed': 358, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Updated best_params for RandomForestRegressor (F2).
$ws.Range("F2").Value = "{'max_depth': 20, 'n_estimators': 100}"

# Updated rmse row (row 3).
$ws.Range("F3").Value = 0.05290137422467665
$ws.Range("G3").Value = 0.05019867482937072
$ws.Range("H3").Value = 0.1015212474423683
$ws.Range("K3").Value = 0.02673389891972606

# Updated r2 row (row 4).
$ws.Range("F4").Value = 0.9483827196779414
$ws.Range("G4").Value = 0.953396796924791
$ws.Range("H4").Value = 0.8230625889627975
$ws.Range("K4").Value = 0.986653353629201

# Updated mape row (row 5).
$ws.Range("F5").Value = 6.427331223082751
$ws.Range("G5").Value = 6.371517571339483
$ws.Range("H5").Value = 19.40265142933292
$ws.Range("K5").Value = 3.774573562720013
